# Update "想去人数" (column F) figures on both the "展览" and "全部类型"
# sheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1159
    6  = 151
    10 = 5291
    11 = 4818
    12 = 17
    13 = 36
    15 = 49
    16 = 188
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
